$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A=20, D=20)
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.312815
$ws.Cells.Item(2, 8).Value = 0.938445
$ws.Cells.Item(2, 9).Value = 0.0082131704949067
$ws.Cells.Item(2, 10).Value = 0.0082131704949067
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.312815
$ws.Cells.Item(2, 14).Value = 0.938445
$ws.Cells.Item(2, 15).Value = 0.0082131704949067
$ws.Cells.Item(2, 16).Value = 0.0082131704949067
$ws.Cells.Item(2, 17).Value = 0.09785322422500001
$ws.Cells.Item(2, 18).Value = 0.8806790180249999
$ws.Cells.Item(2, 19).Value = 0.00006745616957840596
$ws.Cells.Item(2, 20).Value = 0.00006745616957840596

# Row 3 (A=20, D=21)
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.312815
$ws.Cells.Item(3, 8).Value = 0.938445
$ws.Cells.Item(3, 9).Value = 0.0082131704949067
$ws.Cells.Item(3, 10).Value = 0.0082131704949067
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 24.84824866666667
$ws.Cells.Item(3, 14).Value = 74.544746
$ws.Cells.Item(3, 15).Value = 0.6524076620340182
$ws.Cells.Item(3, 16).Value = 0.6524076620340182
$ws.Cells.Item(3, 17).Value = 7.772904906663335
$ws.Cells.Item(3, 18).Value = 69.95614415997
$ws.Cells.Item(3, 19).Value = 0.00535833536046886
$ws.Cells.Item(3, 20).Value = 0.00535833536046886

# Row 4 (A=20, D=22)
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.312815
$ws.Cells.Item(4, 8).Value = 0.938445
$ws.Cells.Item(4, 9).Value = 0.0082131704949067
$ws.Cells.Item(4, 10).Value = 0.0082131704949067
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 12.866992
$ws.Cells.Item(4, 14).Value = 38.600976
$ws.Cells.Item(4, 15).Value = 0.33783162269264755
$ws.Cells.Item(4, 16).Value = 0.3378316226926476
$ws.Cells.Item(4, 17).Value = 4.02498810248
$ws.Cells.Item(4, 18).Value = 36.22489292232
$ws.Cells.Item(4, 19).Value = 0.0027746687157457054
$ws.Cells.Item(4, 20).Value = 0.002774668715745706

# Row 5 (A=20, D=23)
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.312815
$ws.Cells.Item(5, 8).Value = 0.938445
$ws.Cells.Item(5, 9).Value = 0.0082131704949067
$ws.Cells.Item(5, 10).Value = 0.0082131704949067
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.05894133333333334
$ws.Cells.Item(5, 14).Value = 0.176824
$ws.Cells.Item(5, 15).Value = 0.0015475447784274862
$ws.Cells.Item(5, 16).Value = 0.0015475447784274862
$ws.Cells.Item(5, 17).Value = 0.01843773318666667
$ws.Cells.Item(5, 18).Value = 0.16593959868
$ws.Cells.Item(5, 19).Value = 0.000012710249113727556
$ws.Cells.Item(5, 20).Value = 0.000012710249113727556

# Row 6 (A=21, D=20)
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 24.84824866666667
$ws.Cells.Item(6, 8).Value = 74.544746
$ws.Cells.Item(6, 9).Value = 0.6524076620340182
$ws.Cells.Item(6, 10).Value = 0.6524076620340182
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.312815
$ws.Cells.Item(6, 14).Value = 0.938445
$ws.Cells.Item(6, 15).Value = 0.0082131704949067
$ws.Cells.Item(6, 16).Value = 0.0082131704949067
$ws.Cells.Item(6, 17).Value = 7.772904906663335
$ws.Cells.Item(6, 18).Value = 69.95614415997
$ws.Cells.Item(6, 19).Value = 0.00535833536046886
$ws.Cells.Item(6, 20).Value = 0.00535833536046886

# Row 7 (A=21, D=21)
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 24.84824866666667
$ws.Cells.Item(7, 8).Value = 74.544746
$ws.Cells.Item(7, 9).Value = 0.6524076620340182
$ws.Cells.Item(7, 10).Value = 0.6524076620340182
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 24.84824866666667
$ws.Cells.Item(7, 14).Value = 74.544746
$ws.Cells.Item(7, 15).Value = 0.6524076620340182
$ws.Cells.Item(7, 16).Value = 0.6524076620340182
$ws.Cells.Item(7, 17).Value = 617.4354618005019
$ws.Cells.Item(7, 18).Value = 5556.919156204516
$ws.Cells.Item(7, 19).Value = 0.42563575748069377
$ws.Cells.Item(7, 20).Value = 0.42563575748069377

# Row 8 (A=21, D=22)
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 24.84824866666667
$ws.Cells.Item(8, 8).Value = 74.544746
$ws.Cells.Item(8, 9).Value = 0.6524076620340182
$ws.Cells.Item(8, 10).Value = 0.6524076620340182
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 12.866992
$ws.Cells.Item(8, 14).Value = 38.600976
$ws.Cells.Item(8, 15).Value = 0.33783162269264755
$ws.Cells.Item(8, 16).Value = 0.3378316226926476
$ws.Cells.Item(8, 17).Value = 319.7222168080107
$ws.Cells.Item(8, 18).Value = 2877.4999512720965
$ws.Cells.Item(8, 19).Value = 0.22040393912206876
$ws.Cells.Item(8, 20).Value = 0.2204039391220688

# Row 9 (A=21, D=23)
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 24.84824866666667
$ws.Cells.Item(9, 8).Value = 74.544746
$ws.Cells.Item(9, 9).Value = 0.6524076620340182
$ws.Cells.Item(9, 10).Value = 0.6524076620340182
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.05894133333333334
$ws.Cells.Item(9, 14).Value = 0.176824
$ws.Cells.Item(9, 15).Value = 0.0015475447784274862
$ws.Cells.Item(9, 16).Value = 0.0015475447784274862
$ws.Cells.Item(9, 17).Value = 1.464588907411556
$ws.Cells.Item(9, 18).Value = 13.181300166704002
$ws.Cells.Item(9, 19).Value = 0.001009630070786829
$ws.Cells.Item(9, 20).Value = 0.001009630070786829

# Row 10 (A=22, D=20)
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 12.866992
$ws.Cells.Item(10, 8).Value = 38.600976
$ws.Cells.Item(10, 9).Value = 0.33783162269264755
$ws.Cells.Item(10, 10).Value = 0.3378316226926476
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.312815
$ws.Cells.Item(10, 14).Value = 0.938445
$ws.Cells.Item(10, 15).Value = 0.0082131704949067
$ws.Cells.Item(10, 16).Value = 0.0082131704949067
$ws.Cells.Item(10, 17).Value = 4.02498810248
$ws.Cells.Item(10, 18).Value = 36.22489292232
$ws.Cells.Item(10, 19).Value = 0.0027746687157457054
$ws.Cells.Item(10, 20).Value = 0.002774668715745706

# Row 11 (A=22, D=21)
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 12.866992
$ws.Cells.Item(11, 8).Value = 38.600976
$ws.Cells.Item(11, 9).Value = 0.33783162269264755
$ws.Cells.Item(11, 10).Value = 0.3378316226926476
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 24.84824866666667
$ws.Cells.Item(11, 14).Value = 74.544746
$ws.Cells.Item(11, 15).Value = 0.6524076620340182
$ws.Cells.Item(11, 16).Value = 0.6524076620340182
$ws.Cells.Item(11, 17).Value = 319.7222168080107
$ws.Cells.Item(11, 18).Value = 2877.4999512720965
$ws.Cells.Item(11, 19).Value = 0.22040393912206876
$ws.Cells.Item(11, 20).Value = 0.2204039391220688

# Row 12 (A=22, D=22)
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 12.866992
$ws.Cells.Item(12, 8).Value = 38.600976
$ws.Cells.Item(12, 9).Value = 0.33783162269264755
$ws.Cells.Item(12, 10).Value = 0.3378316226926476
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 12.866992
$ws.Cells.Item(12, 14).Value = 38.600976
$ws.Cells.Item(12, 15).Value = 0.33783162269264755
$ws.Cells.Item(12, 16).Value = 0.3378316226926476
$ws.Cells.Item(12, 17).Value = 165.559483128064
$ws.Cells.Item(12, 18).Value = 1490.0353481525763
$ws.Cells.Item(12, 19).Value = 0.11413020529114737
$ws.Cells.Item(12, 20).Value = 0.11413020529114741

# Row 13 (A=22, D=23)
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 12.866992
$ws.Cells.Item(13, 8).Value = 38.600976
$ws.Cells.Item(13, 9).Value = 0.33783162269264755
$ws.Cells.Item(13, 10).Value = 0.3378316226926476
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.05894133333333334
$ws.Cells.Item(13, 14).Value = 0.176824
$ws.Cells.Item(13, 15).Value = 0.0015475447784274862
$ws.Cells.Item(13, 16).Value = 0.0015475447784274862
$ws.Cells.Item(13, 17).Value = 0.7583976644693334
$ws.Cells.Item(13, 18).Value = 6.825578980224001
$ws.Cells.Item(13, 19).Value = 0.0005228095636856914
$ws.Cells.Item(13, 20).Value = 0.0005228095636856915

# Row 14 (A=23, D=20)
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.05894133333333334
$ws.Cells.Item(14, 8).Value = 0.176824
$ws.Cells.Item(14, 9).Value = 0.0015475447784274862
$ws.Cells.Item(14, 10).Value = 0.0015475447784274862
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.312815
$ws.Cells.Item(14, 14).Value = 0.938445
$ws.Cells.Item(14, 15).Value = 0.0082131704949067
$ws.Cells.Item(14, 16).Value = 0.0082131704949067
$ws.Cells.Item(14, 17).Value = 0.01843773318666667
$ws.Cells.Item(14, 18).Value = 0.16593959868
$ws.Cells.Item(14, 19).Value = 0.000012710249113727556
$ws.Cells.Item(14, 20).Value = 0.000012710249113727556

# Row 15 (A=23, D=21)
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.05894133333333334
$ws.Cells.Item(15, 8).Value = 0.176824
$ws.Cells.Item(15, 9).Value = 0.0015475447784274862
$ws.Cells.Item(15, 10).Value = 0.0015475447784274862
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 24.84824866666667
$ws.Cells.Item(15, 14).Value = 74.544746
$ws.Cells.Item(15, 15).Value = 0.6524076620340182
$ws.Cells.Item(15, 16).Value = 0.6524076620340182
$ws.Cells.Item(15, 17).Value = 1.464588907411556
$ws.Cells.Item(15, 18).Value = 13.181300166704002
$ws.Cells.Item(15, 19).Value = 0.001009630070786829
$ws.Cells.Item(15, 20).Value = 0.001009630070786829

# Row 16 (A=23, D=22)
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.05894133333333334
$ws.Cells.Item(16, 8).Value = 0.176824
$ws.Cells.Item(16, 9).Value = 0.0015475447784274862
$ws.Cells.Item(16, 10).Value = 0.0015475447784274862
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 12.866992
$ws.Cells.Item(16, 14).Value = 38.600976
$ws.Cells.Item(16, 15).Value = 0.33783162269264755
$ws.Cells.Item(16, 16).Value = 0.3378316226926476
$ws.Cells.Item(16, 17).Value = 0.7583976644693334
$ws.Cells.Item(16, 18).Value = 6.825578980224001
$ws.Cells.Item(16, 19).Value = 0.0005228095636856914
$ws.Cells.Item(16, 20).Value = 0.0005228095636856915

# Row 17 (A=23, D=23)
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.05894133333333334
$ws.Cells.Item(17, 8).Value = 0.176824
$ws.Cells.Item(17, 9).Value = 0.0015475447784274862
$ws.Cells.Item(17, 10).Value = 0.0015475447784274862
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.05894133333333334
$ws.Cells.Item(17, 14).Value = 0.176824
$ws.Cells.Item(17, 15).Value = 0.0015475447784274862
$ws.Cells.Item(17, 16).Value = 0.0015475447784274862
$ws.Cells.Item(17, 17).Value = 0.0034740807751111116
$ws.Cells.Item(17, 18).Value = 0.031266726976000005
$ws.Cells.Item(17, 19).Value = 0.000002394894841238177
$ws.Cells.Item(17, 20).Value = 0.000002394894841238177
